$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("PUESTOS")

# New column G: header text + same header style ("Accent3") as the rest of row 1
$ws1.Range("G1").Value = "EMAIL USUARIO ASIGNADO"
$ws1.Range("G1").Style = "Accent3"

# Widen the new column to fit its header (closest attainable width to 27.5703125
# given Excel's internal 1/6-character pixel-width snapping)
$ws1.Columns.Item(7).ColumnWidth = 26.666666666666668

# PUESTOS becomes the active sheet/tab, with G5 selected
$ws1.Activate() | Out-Null
$ws1.Range("G5").Select() | Out-Null
